$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("General")

$ids1 = @(
    "C23843_NonMemberCreditCardtLimit",
    "C23889_NonMemberUserNeedToKnowMaxAndMinCreditLimit",
    "C23890_NonMemberUserMustEnterTypeAndLimit",
    "C23838_NonMemberChangeSelectedCreditCardOption",
    "C23761_NonMemberUserToKnowAboutPersonalInfo",
    "C23762_NonMemberUserReceivesPopupForInvalidCharecters"
)

for ($i = 0; $i -lt $ids1.Length; $i++) {
    $row = 48 + $i
    $ws1.Range("A$row").Value = $ids1[$i]
    $ws1.Range("B$row").Value = "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx"
    $ws1.Hyperlinks.Add($ws1.Range("B$row"), "https://forms-preprod.fivision.com/tdecu/oa/Default.aspx") | Out-Null
    $ws1.Range("C$row").Value = "Yes"
    $ws1.Range("D$row").Value = "Chrome"
}
